# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The sheet's "K" column (column G) is being regenerated with newly
# calculated values (s_vals), replacing the previous Strike# derived
# numbers. Write the new computed values for each data row (rows 2-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 4
    3  = 6
    4  = 5
    5  = 13
    6  = 8
    7  = 10
    8  = 10
    9  = 14
    10 = 11
    11 = 8
    12 = 5
    13 = 4
    14 = 8
    15 = 5
    16 = 8
    17 = 11
    18 = 6
    19 = 10
    20 = 6
    21 = 5
    22 = 10
    23 = 13
    24 = 9
    25 = 5
    26 = 7
    27 = 9
    28 = 10
    29 = 12
    30 = 5
    31 = 11
    32 = 3
    33 = 3
    34 = 7
    35 = 3
    36 = 6
    37 = 2
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}
